$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.6514076380946
$ws.Range("C2").Value = 7.741825815864408
$ws.Range("D2").Value = 13.26019950611015
$ws.Range("E2").Value = 13.48439481454401
$ws.Range("G2").Value = 3.705182133722158
$ws.Range("I2").Value = 29.33246696157861
$ws.Range("J2").Value = 8.319089616753658
$ws.Range("L2").Value = 12.63027525088045
$ws.Range("M2").Value = 18.16039445697825
$ws.Range("O2").Value = 31.19983929856758
$ws.Range("B3").Value = 18.24080647313166
$ws.Range("C3").Value = 7.369660582656021
$ws.Range("D3").Value = 13.28222574096133
$ws.Range("E3").Value = 13.52233879294058
$ws.Range("G3").Value = 3.707660028416202
$ws.Range("I3").Value = 29.47480329809872
$ws.Range("J3").Value = 8.316194228781148
$ws.Range("L3").Value = 12.63618920080549
$ws.Range("M3").Value = 18.07064374584845
$ws.Range("O3").Value = 31.31108256915581
$ws.Range("B4").Value = 17.98629201786513
$ws.Range("C4").Value = 7.129972900285219
$ws.Range("D4").Value = 13.29799462284642
$ws.Range("E4").Value = 13.54694938079308
$ws.Range("G4").Value = 3.709262343049923
$ws.Range("I4").Value = 29.56815191760012
$ws.Range("J4").Value = 8.314475424074157
$ws.Range("L4").Value = 12.64138520173457
$ws.Range("M4").Value = 18.01759750333992
$ws.Range("O4").Value = 31.3864662455634
$ws.Range("B5").Value = 17.88211093820032
$ws.Range("C5").Value = 7.029560876870465
$ws.Range("D5").Value = 13.30498483565549
$ws.Range("E5").Value = 13.5573093356443
$ws.Range("G5").Value = 3.709935703182679
$ws.Range("I5").Value = 29.60768834233902
$ws.Range("J5").Value = 8.313789943985206
$ws.Range("L5").Value = 12.64389688025149
$ws.Range("M5").Value = 17.99651377730838
$ws.Range("O5").Value = 31.4189614583131
$ws.Range("B6").Value = 17.86478789424742
$ws.Range("C6").Value = 7.012724602084848
$ws.Range("D6").Value = 13.30617962843176
$ws.Range("E6").Value = 13.55904961035747
$ws.Range("G6").Value = 3.710048748425966
$ws.Range("I6").Value = 29.61434366976518
$ws.Range("J6").Value = 8.313677026708271
$ws.Range("L6").Value = 12.64433777653662
$ws.Range("M6").Value = 17.99304545484168
$ws.Range("O6").Value = 31.42446437988886
$ws.Range("B7").Value = 17.9848886868479
$ws.Range("C7").Value = 7.128629685872015
$ws.Range("D7").Value = 13.29808661091065
$ws.Range("E7").Value = 13.54708775771241
$ws.Range("G7").Value = 3.70927134151693
$ws.Range("I7").Value = 29.56867906324241
$ws.Range("J7").Value = 8.314466118776496
$ws.Range("L7").Value = 12.64141747776911
$ws.Range("M7").Value = 18.01731098350601
$ws.Range("O7").Value = 31.38689730353673
$ws.Range("B8").Value = 18.51041362479497
$ws.Range("C8").Value = 7.615862960906831
$ws.Range("D8").Value = 13.26732818481191
$ws.Range("E8").Value = 13.49720591334867
$ws.Range("G8").Value = 3.706019765926944
$ws.Range("I8").Value = 29.38030828806662
$ws.Range("J8").Value = 8.318079115935513
$ws.Range("L8").Value = 12.63199009372077
$ws.Range("M8").Value = 18.12902968888907
$ws.Range("O8").Value = 31.23672386231295
$ws.Range("B9").Value = 19.51602563130382
$ws.Range("C9").Value = 8.480037912007216
$ws.Range("D9").Value = 13.22482707642414
$ws.Range("E9").Value = 13.40976714634076
$ws.Range("G9").Value = 3.700282120141809
$ws.Range("I9").Value = 29.05819579223994
$ws.Range("J9").Value = 8.325627416090542
$ws.Range("L9").Value = 12.62588287671154
$ws.Range("M9").Value = 18.36378854070457
$ws.Range("O9").Value = 30.99860818481773
$ws.Range("B10").Value = 20.2322389536358
$ws.Range("C10").Value = 9.056525184111115
$ws.Range("D10").Value = 13.20446641917965
$ws.Range("E10").Value = 13.35180113682492
$ws.Range("G10").Value = 3.696451775583361
$ws.Range("I10").Value = 28.85042450892835
$ws.Range("J10").Value = 8.331449827102615
$ws.Range("L10").Value = 12.62889149814105
$ws.Range("M10").Value = 18.54488641102255
$ws.Range("O10").Value = 30.85829659434445
$ws.Range("B11").Value = 20.55172885545459
$ws.Range("C11").Value = 9.305682666409458
$ws.Range("D11").Value = 13.19756189029624
$ws.Range("E11").Value = 13.326782433508
$ws.Range("G11").Value = 3.694791968606451
$ws.Range("I11").Value = 28.76219057249487
$ws.Range("J11").Value = 8.334157479327784
$ws.Range("L11").Value = 12.63187517362371
$ws.Range("M11").Value = 18.62893475811702
$ws.Range("O11").Value = 30.80203928062848
$ws.Range("B12").Value = 20.67169422055946
$ws.Range("C12").Value = 9.398123763987178
$ws.Range("D12").Value = 13.19528604615318
$ws.Range("E12").Value = 13.31750184952246
$ws.Range("G12").Value = 3.694175256995889
$ws.Range("I12").Value = 28.72968358101729
$ws.Range("J12").Value = 8.335191166500827
$ws.Range("L12").Value = 12.63323600732842
$ws.Range("M12").Value = 18.66098289687191
$ws.Range("O12").Value = 30.78182905187188
$ws.Range("B13").Value = 20.64590449708719
$ws.Range("C13").Value = 9.378300216542474
$ws.Range("D13").Value = 13.1957611294432
$ws.Range("E13").Value = 13.31949199669137
$ws.Range("G13").Value = 3.694307552099091
$ws.Range("I13").Value = 28.73664424925064
$ws.Range("J13").Value = 8.334968173781624
$ws.Range("L13").Value = 12.63293267556327
$ws.Range("M13").Value = 18.65407123556188
$ws.Range("O13").Value = 30.7861329967377
$ws.Range("B14").Value = 20.56161949217036
$ws.Range("C14").Value = 9.313326255083682
$ws.Range("D14").Value = 13.1973678681271
$ws.Range("E14").Value = 13.32601504122158
$ws.Range("G14").Value = 3.694740994806804
$ws.Range("I14").Value = 28.75949804828024
$ws.Range("J14").Value = 8.334242354617638
$ws.Range("L14").Value = 12.63198250858186
$ws.Range("M14").Value = 18.63156705985959
$ws.Range("O14").Value = 30.80035463869212
$ws.Range("B15").Value = 20.50985672533146
$ws.Range("C15").Value = 9.273278451184593
$ws.Range("D15").Value = 13.19839614728012
$ws.Range("E15").Value = 13.33003576682423
$ws.Range("G15").Value = 3.695008028636403
$ws.Range("I15").Value = 28.77361461284722
$ws.Range("J15").Value = 8.333798854095411
$ws.Range("L15").Value = 12.63143054538887
$ws.Range("M15").Value = 18.61781080325554
$ws.Range("O15").Value = 30.809208301063
$ws.Range("B16").Value = 20.21122281044593
$ws.Range("C16").Value = 9.039976309819483
$ws.Range("D16").Value = 13.20496506720858
$ws.Range("E16").Value = 13.35346327825651
$ws.Range("G16").Value = 3.696561905436655
$ws.Range("I16").Value = 28.85631738020928
$ws.Range("J16").Value = 8.331274052196116
$ws.Range("L16").Value = 12.62872890102998
$ws.Range("M16").Value = 18.53942558579096
$ws.Range("O16").Value = 30.86212587924499
$ws.Range("B17").Value = 20.02632144198541
$ws.Range("C17").Value = 8.893479303359758
$ws.Range("D17").Value = 13.20959857493498
$ws.Range("E17").Value = 13.3681806340611
$ws.Range("G17").Value = 3.697536279549179
$ws.Range("I17").Value = 28.90866321582489
$ws.Range("J17").Value = 8.329740191117752
$ws.Range("L17").Value = 12.62748431258155
$ws.Range("M17").Value = 18.49175228886281
$ws.Range("O17").Value = 30.89653148818336
$ws.Range("B18").Value = 19.91938234863047
$ws.Range("C18").Value = 8.807987288040209
$ws.Range("D18").Value = 13.21248558145158
$ws.Range("E18").Value = 13.37677280931967
$ws.Range("G18").Value = 3.698104495514244
$ws.Range("I18").Value = 28.93936235311532
$ws.Range("J18").Value = 8.328863516589244
$ws.Range("L18").Value = 12.62692055449925
$ws.Range("M18").Value = 18.46448982552527
$ws.Range("O18").Value = 30.91703305833689
$ws.Range("B19").Value = 19.88307697245342
$ws.Range("C19").Value = 8.778830689426547
$ws.Range("D19").Value = 13.21350119636961
$ws.Range("E19").Value = 13.37970383094688
$ws.Range("G19").Value = 3.698298222086319
$ws.Range("I19").Value = 28.94985803559511
$ws.Range("J19").Value = 8.328567648766017
$ws.Range("L19").Value = 12.62675583526648
$ws.Range("M19").Value = 18.45528691793577
$ws.Range("O19").Value = 30.92409674943949
$ws.Range("B20").Value = 20.04606625895473
$ws.Range("C20").Value = 8.909201758929381
$ws.Range("D20").Value = 13.20908236227227
$ws.Range("E20").Value = 13.36660079298145
$ws.Range("G20").Value = 3.697431750828107
$ws.Range("I20").Value = 28.90302971200505
$ws.Range("J20").Value = 8.32990289854297
$ws.Range("L20").Value = 12.62760106858652
$ws.Range("M20").Value = 18.49681097939081
$ws.Range("O20").Value = 30.89279519132686
$ws.Range("B21").Value = 20.58640452002781
$ws.Range("C21").Value = 9.332462693252431
$ws.Range("D21").Value = 13.19688673916985
$ws.Range("E21").Value = 13.32409382027236
$ws.Range("G21").Value = 3.694613361840042
$ws.Range("I21").Value = 28.7527607487558
$ws.Range("J21").Value = 8.334455319483475
$ws.Range("L21").Value = 12.63225533740921
$ws.Range("M21").Value = 18.63817123857568
$ws.Range("O21").Value = 30.79614769393269
$ws.Range("B22").Value = 20.93356466555061
$ws.Range("C22").Value = 9.597949664345306
$ws.Range("D22").Value = 13.19089048709652
$ws.Range("E22").Value = 13.29744033807675
$ws.Range("G22").Value = 3.692840256806262
$ws.Range("I22").Value = 28.65982893208976
$ws.Range("J22").Value = 8.337479247203976
$ws.Range("L22").Value = 12.63664285300005
$ws.Range("M22").Value = 18.73183675654931
$ws.Range("O22").Value = 30.73935689227978
$ws.Range("B23").Value = 20.74886038051601
$ws.Range("C23").Value = 9.457281034707782
$ws.Range("D23").Value = 13.19391027185004
$ws.Range("E23").Value = 13.31156289223888
$ws.Range("G23").Value = 3.69378031452272
$ws.Range("I23").Value = 28.70894483715
$ws.Range("J23").Value = 8.335860910150082
$ws.Range("L23").Value = 12.63417846345535
$ws.Range("M23").Value = 18.6817349925207
$ws.Range("O23").Value = 30.76908256559052
$ws.Range("B24").Value = 20.03714160235297
$ws.Range("C24").Value = 8.902097587169729
$ws.Range("D24").Value = 13.20931504694032
$ws.Range("E24").Value = 13.36731463111176
$ws.Range("G24").Value = 3.697478983227807
$ws.Range("I24").Value = 28.90557473518144
$ws.Range("J24").Value = 8.329829322425704
$ws.Range("L24").Value = 12.62754781034397
$ws.Range("M24").Value = 18.49452349059397
$ws.Range("O24").Value = 30.89448212452035
$ws.Range("B25").Value = 19.24744274367076
$ws.Range("C25").Value = 8.256349616495005
$ws.Range("D25").Value = 13.23441662641086
$ws.Range("E25").Value = 13.4323159016194
$ws.Range("G25").Value = 3.701766371335616
$ws.Range("I25").Value = 29.14026795534402
$ws.Range("J25").Value = 8.323536397502258
$ws.Range("L25").Value = 12.62621466088102
$ws.Range("M25").Value = 18.29869626668553
$ws.Range("O25").Value = 31.05696165639532
